# chore: update Sheets via scheduled runner
#
# Refreshes the cached crafting-profit figures (columns H:N - the
# NQ/HQ price + profit/loss projections) for a handful of leve rows
# across several job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW), as if
# a scheduled market-board pull recomputed them with newer prices.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H28").Value = 564.1724
$ws.Range("I28").Value = 550.5238000000001
$ws.Range("J28").Value = 600
$ws.Range("K28").Value = 550.5238000000001
$ws.Range("L28").Value = 600
$ws.Range("M28").Value = -65.52380000000005
$ws.Range("N28").Value = -1570

$ws.Range("I43").Value = 1293.5
$ws.Range("J43").Value = 2337.25
$ws.Range("K43").Value = 1293.5
$ws.Range("L43").Value = 2337.25
$ws.Range("M43").Value = -1224.5
$ws.Range("N43").Value = -2475.25

$ws.Range("H55").Value = 205
$ws.Range("J55").Value = 205.55556
$ws.Range("L55").Value = 205.55556
$ws.Range("N55").Value = -633.55556

$ws.Range("H137").Value = 2217.1064
$ws.Range("I137").Value = 1122.2941
$ws.Range("J137").Value = 5080.4614
$ws.Range("K137").Value = 3366.8823
$ws.Range("L137").Value = 15241.3842
$ws.Range("M137").Value = -816.8823000000002
$ws.Range("N137").Value = -20341.3842

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 3099.487
$ws.Range("I74").Value = 2943.9697
$ws.Range("J74").Value = 3954.8333
$ws.Range("K74").Value = 2943.9697
$ws.Range("L74").Value = 3954.8333
$ws.Range("M74").Value = -2069.9697
$ws.Range("N74").Value = -5702.8333

$ws.Range("H77").Value = 3099.487
$ws.Range("I77").Value = 2943.9697
$ws.Range("J77").Value = 3954.8333
$ws.Range("K77").Value = 14719.8485
$ws.Range("L77").Value = 19774.1665
$ws.Range("M77").Value = -10351.8485
$ws.Range("N77").Value = -28510.1665

$ws.Range("H102").Value = 1731.4445
$ws.Range("I102").Value = 1640.4286
$ws.Range("K102").Value = 1640.4286
$ws.Range("M102").Value = -18.42859999999996

$ws.Range("H110").Value = 1328.68
$ws.Range("I110").Value = 1214.4762
$ws.Range("J110").Value = 1928.25
$ws.Range("K110").Value = 1214.4762
$ws.Range("L110").Value = 1928.25
$ws.Range("M110").Value = 830.5237999999999
$ws.Range("N110").Value = -6018.25

$ws.Range("H132").Value = 2349.6155
$ws.Range("I132").Value = 1728.8
$ws.Range("J132").Value = 4419
$ws.Range("K132").Value = 5186.4
$ws.Range("L132").Value = 13257
$ws.Range("M132").Value = -2656.4
$ws.Range("N132").Value = -18317

# --- BSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H105").Value = 1844.3158
$ws.Range("J105").Value = 2020.5
$ws.Range("L105").Value = 2020.5
$ws.Range("N105").Value = -5514.5

$ws.Range("H107").Value = 1835.4667
$ws.Range("I107").Value = 1703.5555
$ws.Range("J107").Value = 2033.3334
$ws.Range("K107").Value = 1703.5555
$ws.Range("L107").Value = 2033.3334
$ws.Range("M107").Value = 216.4445000000001
$ws.Range("N107").Value = -5873.3334

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 9436076
$ws.Range("I31").Value = 1126.6757
$ws.Range("J31").Value = 31254396
$ws.Range("K31").Value = 1126.6757
$ws.Range("L31").Value = 31254396
$ws.Range("M31").Value = -831.6757
$ws.Range("N31").Value = -31254986

$ws.Range("H34").Value = 9436076
$ws.Range("I34").Value = 1126.6757
$ws.Range("J34").Value = 31254396
$ws.Range("K34").Value = 1126.6757
$ws.Range("L34").Value = 31254396
$ws.Range("M34").Value = -924.6757
$ws.Range("N34").Value = -31254800

$ws.Range("H41").Value = 43813
$ws.Range("J41").Value = 43813
$ws.Range("L41").Value = 43813
$ws.Range("N41").Value = -44669

$ws.Range("H50").Value = 25352.7
$ws.Range("J50").Value = 25352.7
$ws.Range("L50").Value = 25352.7
$ws.Range("N50").Value = -26602.7

$ws.Range("H51").Value = 32322.572
$ws.Range("J51").Value = 32322.572
$ws.Range("L51").Value = 32322.572
$ws.Range("N51").Value = -33794.572

$ws.Range("H59").Value = 37626
$ws.Range("J59").Value = 37626
$ws.Range("L59").Value = 37626
$ws.Range("N59").Value = -39916

$ws.Range("H60").Value = 31631.818
$ws.Range("J60").Value = 31631.818
$ws.Range("L60").Value = 31631.818
$ws.Range("N60").Value = -32653.818

$ws.Range("H61").Value = 32322.572
$ws.Range("J61").Value = 32322.572
$ws.Range("L61").Value = 32322.572
$ws.Range("N61").Value = -33018.572

$ws.Range("H62").Value = 3976
$ws.Range("I62").Value = 4395
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 4395
$ws.Range("L62").Value = 2300
$ws.Range("M62").Value = -3771
$ws.Range("N62").Value = -3548

$ws.Range("H65").Value = 3976
$ws.Range("I65").Value = 4395
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 21975
$ws.Range("L65").Value = 11500
$ws.Range("M65").Value = -18855
$ws.Range("N65").Value = -17740

$ws.Range("H68").Value = 99999
$ws.Range("J68").Value = 99999
$ws.Range("L68").Value = 99999
$ws.Range("N68").Value = -101497

$ws.Range("H71").Value = 99999
$ws.Range("J71").Value = 99999
$ws.Range("L71").Value = 299997
$ws.Range("N71").Value = -307485

$ws.Range("H99").Value = 11117051
$ws.Range("I99").Value = 22226546
$ws.Range("J99").Value = 7555.5557
$ws.Range("K99").Value = 22226546
$ws.Range("L99").Value = 7555.5557
$ws.Range("M99").Value = -22225048
$ws.Range("N99").Value = -10551.5557

$ws.Range("H126").Value = 11117051
$ws.Range("I126").Value = 22226546
$ws.Range("J126").Value = 7555.5557
$ws.Range("K126").Value = 66679638
$ws.Range("L126").Value = 22666.6671
$ws.Range("M126").Value = -66677168
$ws.Range("N126").Value = -27606.6671

$ws.Range("H132").Value = 2365.9
$ws.Range("I132").Value = 1489
$ws.Range("K132").Value = 4467
$ws.Range("M132").Value = -1937

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H86").Value = 1962.875
$ws.Range("I86").Value = 601.3333
$ws.Range("J86").Value = 2779.8
$ws.Range("K86").Value = 1803.9999
$ws.Range("L86").Value = 8339.400000000001
$ws.Range("M86").Value = -617.9999
$ws.Range("N86").Value = -10711.4

$ws.Range("H89").Value = 1962.875
$ws.Range("I89").Value = 601.3333
$ws.Range("J89").Value = 2779.8
$ws.Range("K89").Value = 5411.9997
$ws.Range("L89").Value = 25018.2
$ws.Range("M89").Value = 516.0002999999997
$ws.Range("N89").Value = -36874.2

$ws.Range("H131").Value = 5435575.5
$ws.Range("I131").Value = 45454876
$ws.Range("J131").Value = 855.7037
$ws.Range("K131").Value = 136364628
$ws.Range("L131").Value = 2567.1111
$ws.Range("M131").Value = -136359588
$ws.Range("N131").Value = -12647.1111

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H97").Value = 1050
$ws.Range("I97").Value = 1050
$ws.Range("K97").Value = 1050
$ws.Range("M97").Value = -554

$ws.Range("H113").Value = 10492
$ws.Range("I113").Value = 1355.7142
$ws.Range("J113").Value = 23282.8
$ws.Range("K113").Value = 1355.7142
$ws.Range("L113").Value = 23282.8
$ws.Range("M113").Value = 814.2858000000001
$ws.Range("N113").Value = -27622.8

$ws.Range("H132").Value = 3199.5652
$ws.Range("I132").Value = 1847.3334
$ws.Range("J132").Value = 5735
$ws.Range("K132").Value = 5542.0002
$ws.Range("L132").Value = 17205
$ws.Range("M132").Value = -3012.0002
$ws.Range("N132").Value = -22265

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H132").Value = 10468.609
$ws.Range("I132").Value = 12947.667
$ws.Range("J132").Value = 7865.6
$ws.Range("K132").Value = 38843.001
$ws.Range("L132").Value = 23596.8
$ws.Range("M132").Value = -36313.001
$ws.Range("N132").Value = -28656.8
